# Horarium Pricing.xlsx - add "Subscription" breakdown table (rows 26-34)
# and tidy up the now-unused tall formatting on row 19 (A19:C19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 19: the A19:C19 cells only held leftover big-font formatting
#    (no values) - clear them out and let the row height shrink back
#    to the sheet default.
# ---------------------------------------------------------------------
$ws.Range("A19:C19").Clear()
$ws.Rows.Item(19).AutoFit()

# ---------------------------------------------------------------------
# 2. Row 26: merged "Subscription" banner, re-using the centered style
#    already used for the F1:G1 / I1:J1 banners.
# ---------------------------------------------------------------------
$ws.Range("F1:G1").Copy()
$ws.Range("A26:C26").PasteSpecial(-4122)
$ws.Range("A26").Value = "Subscription"
$ws.Range("A26:C26").Merge()

# ---------------------------------------------------------------------
# 3. New "Subscription" payout table, rows 28-34 - mirrors the layout
#    of the existing payout table in rows 10-16.
# ---------------------------------------------------------------------

# Row 28 header: "% of total" / "Rs. of total" / "Excluding costs"
$ws.Range("A10:C10").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)
$ws.Range("B28").Value = "% of total"
$ws.Range("C28").Value = "Rs. of total"
$ws.Range("D28").Value = "Excluding costs"

# Row 29: blank % / Rs cells (same style as row 11/12) + Per Year / expances / Excluding expances headers
$ws.Range("A12:C12").Copy()
$ws.Range("A29:C29").PasteSpecial(-4122)
$ws.Range("G29").Value = "Per Year"
$ws.Range("H29").Value = "expances"
$ws.Range("I29").Value = "Excluding expances"

# Row 30: harsh
$ws.Range("A12:C12").Copy()
$ws.Range("A30:C30").PasteSpecial(-4122)
$ws.Range("A30").Value = "harsh"
$ws.Range("B30").Value = 0.35
$ws.Range("C30").Formula = "=`$G`$30*B30"
$ws.Range("D30").Formula = "=`$I`$30*B30"
$ws.Range("G30").Value = 30000
$ws.Range("H30").Value = 4500
$ws.Range("I30").Formula = "=G30-H30"

# Row 31: dev
$ws.Range("A13:C13").Copy()
$ws.Range("A31:C31").PasteSpecial(-4122)
$ws.Range("A31").Value = "dev"
$ws.Range("B31").Value = 0.4
$ws.Range("C31").Formula = "=`$G`$30*B31"
$ws.Range("D31").Formula = "=`$I`$30*B31"

# Row 32: yogi
$ws.Range("A14:C14").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)
$ws.Range("A32").Value = "yogi"
$ws.Range("B32").Value = 0.25
$ws.Range("C32").Formula = "=`$G`$30*B32"
$ws.Range("D32").Formula = "=`$I`$30*B32"

# Row 33: spacer row (same styling as row 15)
$ws.Range("A15:C15").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)

# Row 34: totals row (same styling as row 16)
$ws.Range("A16:C16").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)
$ws.Range("B34").Formula = "=SUM(B29:B32)"
$ws.Range("C34").Formula = "=SUM(C30:C32)"

$ws.Rows.Item("28:34").AutoFit()

# ---------------------------------------------------------------------
# 4. Update the view so the newly-added table is visible/selected.
# ---------------------------------------------------------------------
$ws.Range("D32").Select()
$excel.ActiveWindow.ScrollRow = 22
